$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 4.074, 4.575, 1.642),
    @(3, 3.331, 3.927, 1.868),
    @(4, 3.447, 3.986, 1.577),
    @(5, 3.057, 3.427, 1.591),
    @(6, 3.085, 3.709, 1.866),
    @(7, 2.888, 3.363, 1.558),
    @(8, 3.232, 3.513, 1.55),
    @(9, 3.016, 3.509, 1.881),
    @(10, 2.768, 3.267, 1.471),
    @(11, 3.44, 3.602, 1.345),
    @(12, 3.178, 3.702, 2.314),
    @(13, 2.469, 2.888, 1.234),
    @(14, 2.67, 2.837, 0.899),
    @(15, 2.694, 2.985, 1.101),
    @(16, 2.283, 2.529, 0.875),
    @(17, 2.283, 2.529, 0.875),
    @(18, 2.694, 2.985, 1.101),
    @(19, 2.67, 2.837, 0.899),
    @(20, 2.469, 2.888, 1.234),
    @(21, 3.178, 3.702, 2.314),
    @(22, 3.44, 3.602, 1.345),
    @(23, 2.768, 3.267, 1.471),
    @(24, 3.016, 3.509, 1.881),
    @(25, 3.232, 3.513, 1.55),
    @(26, 2.888, 3.363, 1.558),
    @(27, 3.085, 3.709, 1.866),
    @(28, 3.057, 3.427, 1.591),
    @(29, 3.447, 3.986, 1.577),
    @(30, 3.331, 3.927, 1.868),
    @(31, 4.074, 4.575, 1.642),
    @(32, 4.549, 4.162, 16.018),
    @(33, 5.004, 3.963, 15.113),
    @(34, 4.434, 4.066, 11.76),
    @(35, 4.346, 4.564, 17.766),
    @(36, 4.145, 4.336, 12.518),
    @(37, 4.149, 4.347, 10.029),
    @(38, 4.339, 4.593, 17.548),
    @(39, 4.211, 4.468, 12.717),
    @(40, 4.149, 4.379, 10.087),
    @(41, 5.089, 5.414, 13.116),
    @(42, 4.253, 4.478, 12.843),
    @(43, 3.649, 3.813, 4.309),
    @(44, 4.258, 4.478, 12.86),
    @(45, 4.258, 4.475, 12.858),
    @(46, 4.258, 4.478, 12.86),
    @(47, 3.649, 3.813, 4.309),
    @(48, 4.253, 4.478, 12.843),
    @(49, 5.089, 5.414, 13.116),
    @(50, 4.149, 4.379, 10.087),
    @(51, 4.211, 4.468, 12.717),
    @(52, 4.339, 4.593, 17.548),
    @(53, 4.149, 4.347, 10.029),
    @(54, 4.145, 4.336, 12.518),
    @(55, 4.346, 4.564, 17.766),
    @(56, 4.434, 4.066, 11.76),
    @(57, 5.004, 3.963, 15.113),
    @(58, 4.549, 4.162, 16.018),
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}

Write-Output "Updated $($data.Count) rows with 3-digit rounded values"
